$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Magnesium chloride unit price: update baseline value and convert
# the Lower/Upper formulas into their (newly recomputed) static values
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 - Zinc sulfate unit price: update baseline value and convert the
# Lower/Upper formulas into their (newly recomputed) static values
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Row 18 - Fermentation magnesium chloride loading: change Lower formula
# multiplier from 0.08 to 0.8
$ws.Range("G18").Formula = "=E18*0.8"

# Row 19 - Fermentation zinc sulfate loading: change Lower formula
# multiplier from 0.08 to 0.8
$ws.Range("G19").Formula = "=E19*0.8"

# Re-fill Q8's helper formula from Q7 so it collapses back into the
# Q5:Q32 shared-formula group (matches Excel's own behavior of merging
# adjacent identical shared formulas back together on save).
$ws.Range("Q7:Q8").FillDown()

# Update last active cell selection to match the authored state
$ws.Range("F8").Select()
